$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.896.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.437.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.14%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +12.15%  '
$ws.Range('E10').Value = '  -1.72%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000177'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '68.786.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.883.95'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.22'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.436.86'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '339.41'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.89%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.70'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.563.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0822'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '428.90'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.11%  '
$ws.Range('E33').Value = '  +2.96%  '
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '160.19'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.02'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.56%  '
$ws.Range('E39').Value = '  -2.30%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.51'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('E43').Value = '  +0.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.94%  '
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '129.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0720'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('E51').Value = '  +0.73%  '
